$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Split "...detect radioactive contamination oneself or one's food or
#    water" into "...contamination " + "on or in " + "oneself or one's
#    food or water" by inserting "on or in " at the boundary.
# ---------------------------------------------------------------------
$target = $d.Content
$target.Find.Execute("radioactive contamination oneself or one's food or water") | Out-Null

$prefix = "radioactive contamination "
$splitPos = $target.Start + $prefix.Length
$tailEnd = $target.End

# Insert the new text at the split point; this naturally creates a new
# run for the inserted text, splitting the original run in two.
$insertion = $d.Range($splitPos, $splitPos)
$insertion.InsertBefore("on or in ")

$insertedLen = "on or in ".Length
$run2 = $d.Range($splitPos, $splitPos + $insertedLen)
$run3 = $d.Range($splitPos + $insertedLen, $tailEnd + $insertedLen)

# Nudge formatting on the two new runs (toggling Bold on/off) so they are
# recorded as their own distinct runs rather than re-merging into the
# surrounding text when the paragraph is re-serialized.
$run2.Font.Bold = 1
$run2.Font.Bold = 0
$run3.Font.Bold = 1
$run3.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) After the paragraph ending " (Ikegami 2012:155)." (the one right
#    before the "Research Proposal" heading) - the same paragraph that
#    starts with "Stratification plays a significant role..." - add
#    three new note paragraphs.
# ---------------------------------------------------------------------
$notesPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Stratification plays a significant role*") {
        $notesPara = $p
        break
    }
}

$notesPara.Range.InsertParagraphAfter()
$p1 = $notesPara.Next()
$p1.Range.Text = "### work as identity; 64-year-old organic cabbage farmer who committed suicide ###"

$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "### not the first exposure to radiation in culture's history: Hiroshima,Nagasaki; hibakusha...###"

$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "### Confucianism and mistrust of government ###"

Write-Host "Done."
